$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number (e.g. "212.35") must keep
# their original text type, so force the Text number format before writing.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "26.193.76"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.602.48"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "212.35"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D7").Value = "0.484"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").Value = "18.28"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "1.825.98"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.607.72"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "0.511"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "26.169.42"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "61.82"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "200.29"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "9.27"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "1.86"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("D25").Value = "144.28"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").Value = "15.17"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "6.55"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").Value = "0.0488"
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").Value = "2.93"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "1.48"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.39"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("D36").Value = "1.162.12"
$ws.Range("E36").Value = "  +4.53%  "
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "0.785"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").Value = "0.495"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "5.30"
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("D44").Value = "1.738.76"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "91.44"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0103"
$ws.Range("E47").Value = "  +15.63%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "54.02"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  -0.09%  "
